$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "Submodule"
$ws.Range("B31").Value = "Submodule basic"
$ws.Range("C31").Value = @'
When there is git under another git folder, the upper git will regard the deeper git as "Submodule". But not until in upper git invoke command: $ git submodule add {deeper_git} then the upper git will generate the .gitmodule file and start the management of the deeper git module.
But the submodule's content will not commit to the upper git anyway. 
Here are some usage:
$ git submodule add /path/to/git/name.git
$ git submodule status
$ git submodule init
$ git submodule deinit        // delete the git repo
'@

$ws.Rows(31).RowHeight = 115.5

$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
